$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 8 (year 2025) metrics with refreshed BIBI data
$ws.Range("C8").Value = 1375
$ws.Range("D8").Value = 213
$ws.Range("E8").Value = 1162
$ws.Range("F8").Value = 8.736669401148482
$ws.Range("G8").Value = 84.50909090909092
$ws.Range("H8").Value = 15.49090909090909
